$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format so numeric-looking strings (e.g. "1.00") are not
# auto-converted to numbers, then reset the style afterwards so no stray
# number-format style is left applied to the cell.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '65.612.67'
Set-TextValue 'E2' '  -0.85%  '
Set-TextValue 'D3' '3.435.56'
Set-TextValue 'E3' '  -3.10%  '
Set-TextValue 'E4' '  +0.01%  '
Set-TextValue 'D5' '590.90'
Set-TextValue 'E5' '  -1.98%  '
Set-TextValue 'D6' '137.71'
Set-TextValue 'E6' '  -6.12%  '
Set-TextValue 'D7' '3.434.38'
Set-TextValue 'E7' '  -3.14%  '
Set-TextValue 'D8' '1.00'
Set-TextValue 'E8' '  +0.06%  '
Set-TextValue 'D9' '0.501'
Set-TextValue 'E9' '  -0.18%  '
Set-TextValue 'E10' '  -6.42%  '
Set-TextValue 'E11' '  -8.50%  '
Set-TextValue 'D12' '0.378'
Set-TextValue 'E12' '  -7.49%  '
Set-TextValue 'D13' '4.016.77'
Set-TextValue 'E13' '  -3.05%  '
Set-TextValue 'E14' '  -10.40%  '
Set-TextValue 'D15' '26.43'
Set-TextValue 'E15' '  -9.06%  '
Set-TextValue 'D16' '3.443.57'
Set-TextValue 'E16' '  -2.58%  '
Set-TextValue 'D17' '65.557.33'
Set-TextValue 'E17' '  -0.86%  '
Set-TextValue 'E18' '  -1.63%  '
Set-TextValue 'D19' '9.79'
Set-TextValue 'E19' '  -11.36%  '
Set-TextValue 'D20' '5.91'
Set-TextValue 'E20' '  -5.42%  '
Set-TextValue 'D21' '13.68'
Set-TextValue 'E21' '  -6.70%  '
Set-TextValue 'D22' '392.88'
Set-TextValue 'E22' '  -6.41%  '
Set-TextValue 'D23' '0.553'
Set-TextValue 'E23' '  -8.30%  '
Set-TextValue 'D24' '73.44'
Set-TextValue 'E24' '  -5.69%  '
Set-TextValue 'E25' '  -0.06%  '
Set-TextValue 'D26' '3.575.49'
Set-TextValue 'E26' '  -2.88%  '
Set-TextValue 'E27' '  -8.13%  '
Set-TextValue 'E28' '  +0.29%  '
Set-TextValue 'D29' '8.27'
Set-TextValue 'E29' '  -9.58%  '
Set-TextValue 'D30' '7.17'
Set-TextValue 'E30' '  -8.49%  '
Set-TextValue 'E31' '  -9.31%  '
Set-TextValue 'D32' '3.440.75'
Set-TextValue 'E32' '  -2.78%  '
Set-TextValue 'E33' '  -0.01%  '
Set-TextValue 'D34' '0.144'
Set-TextValue 'E34' '  -7.04%  '
Set-TextValue 'D35' '23.00'
Set-TextValue 'E35' '  -6.04%  '
Set-TextValue 'D36' '172.77'
Set-TextValue 'E36' '  -0.81%  '
Set-TextValue 'D37' '6.89'
Set-TextValue 'E37' '  -9.37%  '
Set-TextValue 'E38' '  -9.45%  '
Set-TextValue 'E39' '  -8.02%  '
Set-TextValue 'E40' '  -9.19%  '
Set-TextValue 'D41' '0.0761'
Set-TextValue 'E41' '  -7.90%  '
Set-TextValue 'E42' '  -4.33%  '
Set-TextValue 'D43' '43.80'
Set-TextValue 'E43' '  -4.03%  '
Set-TextValue 'D44' '0.999'
Set-TextValue 'E44' '  +0.01%  '
Set-TextValue 'D45' '4.42'
Set-TextValue 'E45' '  -13.49%  '
Set-TextValue 'E46' '  -10.32%  '
Set-TextValue 'D47' '1.11'
Set-TextValue 'E47' '  +1.15%  '
Set-TextValue 'D48' '22.93'
Set-TextValue 'E48' '  +1.01%  '
Set-TextValue 'D49' '6.56'
Set-TextValue 'E49' '  -8.09%  '
Set-TextValue 'E50' '  -12.74%  '
Set-TextValue 'D51' '2.211.88'
Set-TextValue 'E51' '  -6.95%  '
